# Updated symbol list on Mon Dec 12 11:27:44 UTC 2022 with GitHub Actions
#
# This script updates the "Price" column (D) for a number of coins, and
# swaps the BKEXToken / CEJI rows (42 and 43), updating their Coin name,
# Link, Price and Volume(1h) columns accordingly.
#
# Prices in column D are stored as text (not numbers) in the workbook, so
# each numeric-looking value is written with a leading apostrophe to force
# Excel to keep it as text instead of re-interpreting it as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.Value = "'" + $value
}

# --- Price (column D) updates -------------------------------------------
Set-TextValue $ws.Range("D2")  "283.53"
Set-TextValue $ws.Range("D3")  "20.82"
Set-TextValue $ws.Range("D4")  "6.215"
Set-TextValue $ws.Range("D8")  "1.479"
Set-TextValue $ws.Range("D9")  "0.8170"
Set-TextValue $ws.Range("D11") "0.1646"
Set-TextValue $ws.Range("D13") "0.03606"
Set-TextValue $ws.Range("D14") "0.03135"
Set-TextValue $ws.Range("D15") "0.09143"
Set-TextValue $ws.Range("D17") "0.001642"
Set-TextValue $ws.Range("D18") "0.04665"
Set-TextValue $ws.Range("D19") "0.006431"
Set-TextValue $ws.Range("D24") "2.323"
Set-TextValue $ws.Range("D40") "0.04714"
Set-TextValue $ws.Range("D41") "0.007048"
Set-TextValue $ws.Range("D45") "0.00006362"
Set-TextValue $ws.Range("D47") "0.9995"
Set-TextValue $ws.Range("D48") "0.002729"

# --- Rows 42/43: BKEXToken and CEJI swap places --------------------------
# Row 42 becomes CEJI (previously BKEXToken)
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D42") "0.004398"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 becomes BKEXToken (previously CEJI)
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D43") "0.1104"
$ws.Range("E43").Value = "42BKEXTokenBKK"
